# Insert a new data row at row 58 (pushing existing rows 58-100 down to 59-101)
# and populate it with the new "Arveja Verde" record for Región Metropolitana.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above the current row 58; this shifts rows 58:100 -> 59:101
# and extends the used range to row 101, matching the target dimension A1:R101.
$ws.Rows("58:58").Insert()

# Populate the newly inserted (blank) row 58 with the new record's values.
# (Single-quoted literals throughout so no PowerShell variable interpolation
# can touch the literal '$' in the unit-of-sale text.)
$ws.Range("A58").Value = 10
$ws.Range("B58").Value = 'Vega Modelo de Temuco'
$ws.Range("C58").Value = 'La Araucanía'
$ws.Range("D58").Value = 44879
$ws.Range("E58").Value = 9
$ws.Range("F58").Value = 100112022
$ws.Range("G58").Value = 'Arveja Verde'
$ws.Range("H58").Value = 'Sin especificar'
$ws.Range("I58").Value = 'Primera'
$ws.Range("J58").Value = 160
$ws.Range("K58").Value = 22000
$ws.Range("L58").Value = 23000
$ws.Range("M58").Value = 22500
$ws.Range("N58").Value = '$/saco 25 kilos'
$ws.Range("O58").Value = 'Región Metropolitana'
$ws.Range("P58").Value = 900
$ws.Range("Q58").Value = 25
$ws.Range("R58").Value = 'Hortaliza'
